# Actualización desde MV -datos-
# Appends 14 new daily rows (17-09-2021 .. 30-09-2021) below the existing
# data, following the same pattern as the prior rows (column B = 449,
# column C = 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startDate = Get-Date -Year 2021 -Month 9 -Day 17
$lastRow = 233

for ($i = 0; $i -lt 14; $i++) {
    $row = $lastRow + 1 + $i
    $date = $startDate.AddDays($i)
    $dateText = $date.ToString("dd-MM-yyyy")

    $ws.Cells.Item($row, 1).Value = $dateText
    $ws.Cells.Item($row, 2).Value = 449
    $ws.Cells.Item($row, 3).Value = 0
}
